# Update the date header (2025-04-13 Sunday -> 2025-04-14 Monday)
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-13 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-14 Monday", 2)

# Update the multiplication answers. Each left-hand string is unique in the
# document except "74×78=5772", which occurs twice (row 15, columns 3 and 5
# of the table) and is handled separately below via direct cell access so
# each occurrence gets its own distinct replacement.

$replacements = @(
    @("16×23=368", "15×24=360"),
    @("85×85=7225", "31×89=2759"),
    @("41×62=2542", "85×11=935"),
    @("95×51=4845", "55×28=1540"),
    @("34×84=2856", "63×35=2205"),
    @("51×55=2805", "21×13=273"),
    @("40×51=2040", "57×99=5643"),
    @("48×17=816", "38×32=1216"),
    @("97×85=8245", "92×84=7728"),
    @("19×68=1292", "54×73=3942"),
    @("73×66=4818", "52×30=1560"),
    @("41×97=3977", "65×19=1235"),
    @("20×94=1880", "17×59=1003"),
    @("93×60=5580", "40×18=720"),
    @("53×52=2756", "21×70=1470"),
    @("41×20=820", "42×77=3234"),
    @("46×18=828", "76×85=6460"),
    @("49×68=3332", "31×58=1798"),
    @("90×23=2070", "79×15=1185"),
    @("52×57=2964", "58×27=1566"),
    @("84×65=5460", "70×33=2310"),
    @("99×85=8415", "29×73=2117"),
    @("11×76=836", "88×31=2728")
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)
}

# The two cells that both originally read "74×78=5772" (table row 15,
# columns 3 and 5, 1-indexed) need different replacement values, so set
# each cell's range text directly instead of using a document-wide Find.
$table = $d.Tables.Item(1)
$table.Cell(15, 3).Range.Text = "18×95=1710"
$table.Cell(15, 5).Range.Text = "41×93=3813"
